# Updates cryptos list data (price + 1h volume change) per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Price/volume cells are stored as plain text (not numbers) in the sheet,
# e.g. "64.475.74" or "1.00" must stay text so trailing zeros and the
# thousands-dot formatting are preserved exactly. Excel auto-converts any
# value that parses as a number, so those are written with a leading
# apostrophe (forces text) and the cell's style is put back to "Normal"
# afterwards so no stray number-format/style is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.475.74"
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "3.661.05"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "'402.04"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").Value = "'125.87"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").Value = "3.654.37"
$ws.Range("E7").Value = "  -4.91%  "
$ws.Range("D8").Value = "'0.594"
$ws.Range("E8").Value = "  -5.63%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'0.699"
$ws.Range("E10").Value = "  -6.49%  "
$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "  -8.25%  "
$ws.Range("D12").Value = "'0.0000338"
$ws.Range("E12").Value = "  -4.93%  "
$ws.Range("D13").Value = "'39.69"
$ws.Range("E13").Value = "  -4.39%  "
$ws.Range("D14").Value = "4.211.79"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "'9.39"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("D16").Value = "'14.48"
$ws.Range("E16").Value = "  +14.05%  "
$ws.Range("D17").Value = "'0.137"
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").Value = "3.651.04"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").Value = "'19.05"
$ws.Range("E19").Value = "  -5.66%  "
$ws.Range("D20").Value = "64.650.02"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "'1.04"
$ws.Range("E21").Value = "  -5.66%  "
$ws.Range("D22").Value = "'400.43"
$ws.Range("E22").Value = "  -8.91%  "
$ws.Range("D23").Value = "'14.20"
$ws.Range("E23").Value = "  -5.14%  "
$ws.Range("D24").Value = "'83.22"
$ws.Range("E24").Value = "  -4.41%  "
$ws.Range("D25").Value = "'2.99"
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'35.29"
$ws.Range("E26").Value = "  -4.23%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.46"
$ws.Range("E27").Value = "  +10.31%  "
$ws.Range("D28").Value = "'3.02"
$ws.Range("E28").Value = "  -9.25%  "
$ws.Range("D29").Value = "'8.94"
$ws.Range("E29").Value = "  -12.00%  "
$ws.Range("D30").Value = "'12.25"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("D32").Value = "'0.114"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("D33").Value = "'6.96"
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("D34").Value = "'0.152"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'37.67"
$ws.Range("E35").Value = "  -6.93%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "'54.78"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").Value = "'0.0444"
$ws.Range("E38").Value = "  -7.31%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0676"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "'2.72"
$ws.Range("E41").Value = "  -8.62%  "
$ws.Range("D42").Value = "'0.132"
$ws.Range("E42").Value = "  -7.64%  "
$ws.Range("D43").Value = "'3.12"
$ws.Range("E43").Value = "  +18.27%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'26.44"
$ws.Range("E44").Value = "  -14.55%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'143.17"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'2.00"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "'3.14"
$ws.Range("E47").Value = "  -6.39%  "
$ws.Range("D48").Value = "'4.15"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "'2.51"
$ws.Range("E49").Value = "  -5.81%  "
$ws.Range("D50").Value = "'2.70"
$ws.Range("E50").Value = "  -5.34%  "
$ws.Range("D51").Value = "'0.284"
$ws.Range("E51").Value = "  -5.52%  "

# Strip the text-format style picked up from the leading apostrophe above
# so these cells end up styled identically to how they started (no <s> diff).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
